# LOQ4236.xlsx content update
# Rebuilds rows 13-22 of the single worksheet with the corrected / expanded
# syllabus content, inserts a new row 22 (Bibliografia) and realigns the
# column A definition.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 13: "Docentes responsáveis:" value row -> A13 becomes blank,
# B13/C13 hold the professor's name.
# ---------------------------------------------------------------------
$ws.Range("A13").Clear()
$ws.Range("B13").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C13").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Rows.Item(13).AutoFit()

# ---------------------------------------------------------------------
# Row 14: "Programa resumido:" + short syllabus text
# ---------------------------------------------------------------------
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "Tópicos que abordem o tema do projeto de seu planejamento a execução."
$ws.Range("C14").Value = "Tópicos que abordem o tema do projeto de seu planejamento a execução."
$ws.Rows.Item(14).RowHeight = 60

# ---------------------------------------------------------------------
# Row 15: "Short syllabus:" label only (no B/C content any more)
# ---------------------------------------------------------------------
$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15:C15").Clear()
$ws.Rows.Item(15).RowHeight = 60

# ---------------------------------------------------------------------
# Row 16: "Programa:" + full program text
# ---------------------------------------------------------------------
$programa = @"
Noções de Gestão de Projetos
Organização do tempo: dimensão pessoal;
Técnicas para a realização de apresentações;
Noções de Aprendizagem Baseada em Projetos
Trabalho em Grupo, Equipes e times. 
Postura e Ética Profissional
Técnicas para redação de relatório técnico;
Tutoria de projetos.
Assuntos Técnicos específicos relacionados com o tema do projeto.
"@
$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa
$ws.Rows.Item(16).RowHeight = 120

# ---------------------------------------------------------------------
# Row 17: "Syllabus:" label only
# ---------------------------------------------------------------------
$ws.Range("A17").Value = "Syllabus:"
$ws.Range("B17:C17").Clear()
$ws.Rows.Item(17).RowHeight = 120

# ---------------------------------------------------------------------
# Row 18: "Avaliação:" label only (loses its old B/C content)
# ---------------------------------------------------------------------
$ws.Range("A18").Value = "Avaliação:"
$ws.Range("B18:C18").Clear()
$ws.Rows.Item(18).AutoFit()

# ---------------------------------------------------------------------
# Row 19: "Método:" + method description (unchanged body text, new index)
# ---------------------------------------------------------------------
$metodo = @"
O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras.

Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Produção, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão. 
Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.
As aulas ocorrerão: 1) através de uma reunião da equipe de trabalho para tratar do projeto, e  2) palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores  ou profissionais de empresas.
"@
$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo
$ws.Rows.Item(19).RowHeight = 60

# ---------------------------------------------------------------------
# Row 20: "Critério:" + grading criteria (unchanged body text, new index)
# ---------------------------------------------------------------------
$criterio = @"
A nota será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros.
O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na coordenação da disciplina.
"@
$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio
$ws.Rows.Item(20).RowHeight = 60

# ---------------------------------------------------------------------
# Row 21: "Norma de recuperação:" + "Não há recuperação" (height 120 -> 60)
# ---------------------------------------------------------------------
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = "Não há recuperação"
$ws.Range("C21").Value = "Não há recuperação"
$ws.Rows.Item(21).RowHeight = 60

# ---------------------------------------------------------------------
# Row 22 (new): "Bibliografia:" + bibliography text
# ---------------------------------------------------------------------
$bibliografia = @"
Artigos sobre metodologias ativas de aprendizagem e  Project Based Learning.
Livros e Artigos científicos relacionados com o tema do projeto.
"@
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
$ws.Rows.Item(22).RowHeight = 120

# ---------------------------------------------------------------------
# Column layout: split the old A:B (1-2) width definition so column A
# alone keeps 30.7109375 and column B keeps its own 60.7109375 entry.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 60.7109375
$ws.Columns.Item(1).ColumnWidth = 30.7109375
